# The edit swaps the full data content of row 3 and row 4 (the two
# observation records in this species-occurrence export) while leaving
# the header row (1) and the first data row (2) untouched.
#
# Rather than relying on a risky whole-range Value array swap (which can
# make Excel "helpfully" reinterpret text such as "1" or "2013-06-03" as
# a number/date), each target cell is written explicitly with its final
# type. Cells that hold text which could be misread as a number or date
# are forced to Text format first so they stay literal strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for row 3 after the edit (this is what row 4 contained
# before the edit).
$row3 = @(
  @("A", "n", "112144586"),
  @("B", "n", "44328"),
  @("C", "s", "Ovaliderad"),
  @("D", "s", "NT"),
  @("E", "n", "102366"),
  @("F", "s", "Ängsmetallvinge"),
  @("G", "s", "Adscita statices"),
  @("H", "s", "(Linnaeus, 1758)"),
  @("P", "s", "Vä, delomr 10, 450 m ONO fängelset, Sk"),
  @("Q", "n", "442995"),
  @("R", "n", "6204827"),
  @("S", "n", "10"),
  @("T", "s", "Skåne"),
  @("U", "s", "Kristianstad"),
  @("V", "s", "Skåne"),
  @("W", "s", "Vä"),
  @("Y", "s", "2013-06-03"),
  @("AA", "s", "2013-06-03"),
  @("AD", "b", "0"),
  @("AE", "b", "0"),
  @("AG", "b", "0"),
  @("AI", "s", "på igenväxande grässandmark"),
  @("AW", "s", "Nils Otto Nilsson"),
  @("AX", "s", "Nils Otto Nilsson"),
  @("AY", "s", "Krst NV-program 2013")
)

# Target state for row 4 after the edit (this is what row 3 contained
# before the edit).
$row4 = @(
  @("A", "n", "112144581"),
  @("B", "n", "42600"),
  @("C", "s", "Ovaliderad"),
  @("D", "s", "NT"),
  @("E", "n", "101260"),
  @("F", "s", "Svartfläckig blåvinge"),
  @("G", "s", "Phengaris arion"),
  @("H", "s", "(Linnaeus, 1758)"),
  @("I", "s", "1"),
  @("J", "s", "ex."),
  @("K", "s", "imago/adult"),
  @("L", "s", "hona"),
  @("M", "s", "vilande"),
  @("P", "s", "Vä, delomr 10, 450 m ONO fängelset, Sk"),
  @("Q", "n", "442972"),
  @("R", "n", "6204767"),
  @("S", "n", "10"),
  @("T", "s", "Skåne"),
  @("U", "s", "Kristianstad"),
  @("V", "s", "Skåne"),
  @("W", "s", "Vä"),
  @("Y", "s", "2013-06-03"),
  @("AA", "s", "2013-06-03"),
  @("AC", "s", "lufthåvning"),
  @("AD", "b", "0"),
  @("AE", "b", "0"),
  @("AG", "b", "0"),
  @("AI", "s", "på igenväxande grässandmark"),
  @("AO", "s", "på grässtrå"),
  @("AQ", "s", "Nils Otto Nilsson"),
  @("AR", "s", "NON 04616"),
  @("AW", "s", "Nils Otto Nilsson"),
  @("AX", "s", "Nils Otto Nilsson"),
  @("AY", "s", "Krst NV-program 2013")
)

function Set-RowData($rowNum, $cellDefs) {
    $fullRange = $ws.Range("A" + $rowNum + ":AY" + $rowNum)
    # Wipe the row first so columns that must end up empty (e.g. those
    # present before the edit but absent afterwards) don't keep stale data.
    $fullRange.ClearContents()

    foreach ($def in $cellDefs) {
        $col = $def[0]
        $kind = $def[1]
        $val = $def[2]
        $cell = $ws.Range($col + $rowNum)

        if ($kind -eq "b") {
            if ($val -eq "1") {
                $cell.Value = $true
            } else {
                $cell.Value = $false
            }
        } elseif ($kind -eq "n") {
            $cell.Value = [double]$val
        } else {
            # Text. Values that look like a plain number (e.g. "1") or a
            # date (e.g. "2013-06-03") would otherwise be silently
            # reinterpreted by Excel as a number/date, so force Text
            # format first for just those risky values to keep them as
            # literal strings.
            if ($val -match '^[0-9]+(\.[0-9]+)?$' -or $val -match '^[0-9]{4}-[0-9]{2}-[0-9]{2}$') {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $val
        }
    }
}

Set-RowData 3 $row3
Set-RowData 4 $row4
